$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("M5").Value = 1.03
$ws.Range("O5").Value = 1.25
$ws.Range("Q5").Value = 1.92
$ws.Range("R5").Value = 1.82
$ws.Range("G6").Value = 1.39
$ws.Range("H6").Value = 4.45
$ws.Range("I6").Value = 6.1
$ws.Range("J6").Value = 1.85
$ws.Range("K6").Value = 2.42
$ws.Range("L6").Value = 5.7
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 2.18
$ws.Range("W6").Value = 7
$ws.Range("X6").Value = 6.3
$ws.Range("Z6").Value = 8
$ws.Range("AC6").Value = 14.5
$ws.Range("AD6").Value = 7.9
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 32
$ws.Range("AJ6").Value = 16
$ws.Range("AK6").Value = 90
$ws.Range("AL6").Value = 45
$ws.Range("AM6").Value = 40
$ws.Range("AT6").Value = 3.25
$ws.Range("AU6").Value = 7.8
$ws.Range("AV6").Value = 65
$ws.Range("AW6").Value = 7.7
$ws.Range("U8").Value = 1.63
$ws.Range("G9").Value = 1.45
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 1.95
$ws.Range("K9").Value = 2.5
$ws.Range("L9").Value = 6
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("Q9").Value = 1.62
$ws.Range("R9").Value = 2.2
$ws.Range("U9").Value = 1.77
$ws.Range("V9").Value = 1.87
$ws.Range("X9").Value = 7.5
$ws.Range("Z9").Value = 10
$ws.Range("AE9").Value = 17
$ws.Range("AJ9").Value = 21
$ws.Range("AN9").Value = 3.5
$ws.Range("AO9").Value = 7
$ws.Range("AU9").Value = 8.5
$ws.Range("AW9").Value = 8
$ws.Range("AX9").Value = 34
$ws.Range("BA9").Value = 126
$ws.Range("AT10").Value = 2.62
$ws.Range("G13").Value = 1.42
